# Natmi following Dr Hou advice
# Update existing rows 2-5 (new sending/target cluster "ECs" added upstream,
# shifting cluster assignments and recalculated edge-weight statistics) and
# append new rows 6-7 for the additional cluster combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntf3"
$ws.Range("C2").Value = "Ntrk3"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 4.477828666666666
$ws.Range("H2").Value = 13.433486
$ws.Range("I2").Value = 0.4652827882180238
$ws.Range("J2").Value = 0.4652827882180238
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1511493333333333
$ws.Range("N2").Value = 0.453448
$ws.Range("O2").Value = 0.7495144539818079
$ws.Range("P2").Value = 0.7495144539818078
$ws.Range("Q2").Value = 0.6768208177475555
$ws.Range("R2").Value = 6.091387359727999
$ws.Range("S2").Value = 0.3487361749583653
$ws.Range("T2").Value = 0.3487361749583652

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntf3"
$ws.Range("C3").Value = "Ntrk3"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 4.477828666666666
$ws.Range("H3").Value = 13.433486
$ws.Range("I3").Value = 0.4652827882180238
$ws.Range("J3").Value = 0.4652827882180238
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.05051366666666667
$ws.Range("N3").Value = 0.151541
$ws.Range("O3").Value = 0.2504855460181921
$ws.Range("P3").Value = 0.2504855460181921
$ws.Range("Q3").Value = 0.2261915446584445
$ws.Range("R3").Value = 2.035723901926
$ws.Range("S3").Value = 0.1165466132596585
$ws.Range("T3").Value = 0.1165466132596585

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ntf3"
$ws.Range("C4").Value = "Ntrk3"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.528563666666667
$ws.Range("H4").Value = 7.585691000000001
$ws.Range("I4").Value = 0.2627383137214249
$ws.Range("J4").Value = 0.2627383137214249
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1511493333333333
$ws.Range("N4").Value = 0.453448
$ws.Range("O4").Value = 0.7495144539818079
$ws.Range("P4").Value = 0.7495144539818078
$ws.Range("Q4").Value = 0.3821907125075556
$ws.Range("R4").Value = 3.439716412568
$ws.Range("S4").Value = 0.1969261637490148
$ws.Range("T4").Value = 0.1969261637490147

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ntf3"
$ws.Range("C5").Value = "Ntrk3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.528563666666667
$ws.Range("H5").Value = 7.585691000000001
$ws.Range("I5").Value = 0.2627383137214249
$ws.Range("J5").Value = 0.2627383137214249
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.05051366666666667
$ws.Range("N5").Value = 0.151541
$ws.Range("O5").Value = 0.2504855460181921
$ws.Range("P5").Value = 0.2504855460181921
$ws.Range("Q5").Value = 0.1277270222034445
$ws.Range("R5").Value = 1.149543199831
$ws.Range("S5").Value = 0.06581214997241017
$ws.Range("T5").Value = 0.06581214997241015

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Ntf3"
$ws.Range("C6").Value = "Ntrk3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.617494
$ws.Range("H6").Value = 7.852482
$ws.Range("I6").Value = 0.2719788980605514
$ws.Range("J6").Value = 0.2719788980605514
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1511493333333333
$ws.Range("N6").Value = 0.453448
$ws.Range("O6").Value = 0.7495144539818079
$ws.Range("P6").Value = 0.7495144539818078
$ws.Range("Q6").Value = 0.395632473104
$ws.Range("R6").Value = 3.560692257936
$ws.Range("S6").Value = 0.203852115274428
$ws.Range("T6").Value = 0.2038521152744279

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Ntf3"
$ws.Range("C7").Value = "Ntrk3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.617494
$ws.Range("H7").Value = 7.852482
$ws.Range("I7").Value = 0.2719788980605514
$ws.Range("J7").Value = 0.2719788980605514
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.05051366666666667
$ws.Range("N7").Value = 0.151541
$ws.Range("O7").Value = 0.2504855460181921
$ws.Range("P7").Value = 0.2504855460181921
$ws.Range("Q7").Value = 0.132219219418
$ws.Range("R7").Value = 1.189972974762
$ws.Range("S7").Value = 0.06812678278612341
$ws.Range("T7").Value = 0.06812678278612341

$wb.Save()
